$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new day column (12-nov) before column DO ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before DO; everything from DO..ES shifts right to DP..ET,
# and the new DO column inherits formatting from its neighbours.
$wsSpot.Range("DO1").EntireColumn.Insert()

# Header for the freshly inserted column.
$wsSpot.Range("DO1").Value = "12-nov"

# Data rows for the new column: no data available yet, so "-" like the other
# not-yet-populated days.
$wsSpot.Range("DO2:DO25").Value = "-"

# --- Sheet "Gaz": append the new daily price row ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A148").NumberFormat = "@"
$wsGaz.Range("A148").Value = "2025-11-10"
$wsGaz.Range("A148").Style = "Normal"
$wsGaz.Range("B148").Value = 28.925

# --- Sheet "CO2": append the new daily price row ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A148").NumberFormat = "@"
$wsCO2.Range("A148").Value = "2025-11-10"
$wsCO2.Range("A148").Style = "Normal"
$wsCO2.Range("B148").Value = 79.88
